# New crime data collected - weekly CompStat update for the 123rd Precinct.
# Bumps the report Volume/Number and the covered week's date range, then
# refreshes the Crime Complaints table (rows 16-30) with the newly
# collected weekly/28-day/YTD/2-year counts and their computed %Chg values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: "Volume 30   Number  20" -> "...Number  22"
# and the covered week "5/15/2023 ... 5/21/2023" -> "5/29/2023 ... 6/4/2023"
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  22"
$ws.Range("C9").Value = "Report Covering the Week  5/29/2023  Through  6/4/2023"

# ---------------------------------------------------------------------
# Cells that flip between the numeric style and the "no data" text
# placeholders ("0" / "***.*"). Copy an already-styled donor cell (which
# keeps the same placeholder state before and after this edit) over the
# target first so the style index lines up exactly, then overwrite the
# value where the target should actually hold a live number.
# ---------------------------------------------------------------------

# C16: was the "0" dash placeholder -> becomes numeric 1 (style like F16)
$ws.Range("F16").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 1

# C18: was the "0" dash placeholder -> becomes numeric 3 (style like D18)
$ws.Range("D18").Copy($ws.Range("C18"))
$ws.Range("C18").Value = 3

# C27: was numeric 1 -> becomes the "0" dash placeholder (style like D27)
$ws.Range("D27").Copy($ws.Range("C27"))

# G27: was numeric 1 -> becomes the "0" dash placeholder (style like D14)
$ws.Range("D14").Copy($ws.Range("G27"))

# H27: was numeric 100 -> becomes the "***.*" placeholder (style like E14)
$ws.Range("E14").Copy($ws.Range("H27"))

# G30: was numeric 1 -> becomes the "0" dash placeholder (style like D14)
$ws.Range("D14").Copy($ws.Range("G30"))

# H30: was numeric -100 -> becomes the "***.*" placeholder (style like E14)
$ws.Range("E14").Copy($ws.Range("H30"))

# ---------------------------------------------------------------------
# Plain numeric refreshes (counts + recomputed %Chg) across the table.
# ---------------------------------------------------------------------

# Row 16 - Robbery
$ws.Range("F16").Value = 2
$ws.Range("I16").Value = 6
$ws.Range("K16").Value = 50
$ws.Range("L16").Value = 100
$ws.Range("M16").Value = -45.454545454545
$ws.Range("N16").Value = -72.727272727272

# Row 17 - Fel. Assault
$ws.Range("F17").Value = 8
$ws.Range("H17").Value = 300
$ws.Range("I17").Value = 36
$ws.Range("J17").Value = 16
$ws.Range("K17").Value = 125
$ws.Range("L17").Value = 125
$ws.Range("M17").Value = 89.473684210526
$ws.Range("N17").Value = -16.279069767441

# Row 18 - Burglary
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 24
$ws.Range("J18").Value = 10
$ws.Range("K18").Value = 140
$ws.Range("L18").Value = 71.428571428571
$ws.Range("M18").Value = -52.941176470588
$ws.Range("N18").Value = -80.487804878048

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -57.142857142857
$ws.Range("F19").Value = 20
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = -23.076923076923
$ws.Range("I19").Value = 123
$ws.Range("J19").Value = 118
$ws.Range("K19").Value = 4.237288135593
$ws.Range("L19").Value = 89.230769230769
$ws.Range("M19").Value = 112.068965517241
$ws.Range("N19").Value = 70.833333333333

# Row 20 - G.L.A.
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 35
$ws.Range("J20").Value = 43
$ws.Range("K20").Value = -18.604651162790
$ws.Range("L20").Value = 133.333333333333
$ws.Range("M20").Value = 133.333333333333
$ws.Range("N20").Value = -88.333333333333

# Row 21 - TOTAL
$ws.Range("C21").Value = 10
$ws.Range("D21").Value = 11
$ws.Range("E21").Value = -9.090909090909
$ws.Range("F21").Value = 39
$ws.Range("H21").Value = 11.428571428571
$ws.Range("I21").Value = 225
$ws.Range("J21").Value = 191
$ws.Range("K21").Value = 17.801047120418
$ws.Range("L21").Value = 97.368421052631
$ws.Range("M21").Value = 44.230769230769
$ws.Range("N21").Value = -60.035523978685

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 9
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 42
$ws.Range("G24").Value = 35
$ws.Range("H24").Value = 20
$ws.Range("I24").Value = 206
$ws.Range("J24").Value = 172
$ws.Range("K24").Value = 19.767441860465
$ws.Range("L24").Value = 110.204081632653
$ws.Range("M24").Value = -2.369668246445

# Row 25 - Misd. Assault
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 12
$ws.Range("G25").Value = 8
$ws.Range("I25").Value = 78
$ws.Range("J25").Value = 75
$ws.Range("K25").Value = 4
$ws.Range("L25").Value = 36.842105263157
$ws.Range("M25").Value = -9.302325581395

# Row 27 - Other Sex Crimes
$ws.Range("F27").Value = 4
$ws.Range("I27").Value = 7
$ws.Range("K27").Value = -22.222222222222
$ws.Range("L27").Value = 40
